$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Haiti / Uganda order (sharedStrings reorder: Uganda now precedes Haiti) ---
# Row 157 previously showed Haiti's data, now shows Uganda's (row moved up in the source order).
# Row 158 previously showed Uganda's data, now shows Haiti's.
$ws.Range("A157").Value = "Uganda"
$ws.Range("A158").Value = "Haiti"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1156744
$ws.Range("C4").Value = 25714
$ws.Range("D4").Value = 160543
$ws.Range("E4").Value = 928973
$ws.Range("G4").Value = 1475
$ws.Range("H4").Value = 67228

# --- Row 9: Alemania ---
$ws.Range("E9").Value = 29155
$ws.Range("G9").Value = 76
$ws.Range("H9").Value = 6812

# --- Row 81: Guinea ---
$ws.Range("B81").Value = 1586
$ws.Range("C81").Value = 49
$ws.Range("D81").Value = 405
$ws.Range("E81").Value = 1174

# --- Row 87: Costa de Marfil ---
$ws.Range("B87").Value = 1362
$ws.Range("C87").Value = 29
$ws.Range("D87").Value = 622
$ws.Range("E87").Value = 725

# --- Row 92: Tunez ---
$ws.Range("B92").Value = 1009
$ws.Range("C92").Value = 11
$ws.Range("D92").Value = 323
$ws.Range("E92").Value = 644
$ws.Range("F92").Value = 25
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 42

# --- Row 157: now Uganda (previously Haiti's slot) ---
$ws.Range("B157").Value = 88
$ws.Range("C157").Value = 3
$ws.Range("D157").Value = 52
$ws.Range("E157").Value = 36
$ws.Range("H157").Value = 0

# --- Row 158: now Haiti (previously Uganda's slot) ---
$ws.Range("D158").Value = 10
$ws.Range("E158").Value = 67
$ws.Range("H158").Value = 8

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Mayo de 2020 a las 23:41"
